$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.965.94"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "'3.903.86"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'466.83"
$ws.Range("E5").Value = "  +9.53%  "
$ws.Range("D6").Value = "'145.44"
$ws.Range("E6").Value = "  +5.18%  "
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.739"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("D11").Value = "'0.0000342"
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("D12").Value = "'43.26"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'10.44"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").Value = "'4.526.89"
$ws.Range("E14").Value = "  +3.32%  "
$ws.Range("D15").Value = "'15.02"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "'3.934.93"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "'20.00"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Value = "'67.221.15"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "'432.86"
$ws.Range("E21").Value = "  +5.02%  "
$ws.Range("D22").Value = "'14.73"
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").Value = "'3.35"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "'88.62"
$ws.Range("E24").Value = "  +3.36%  "
$ws.Range("D25").Value = "'38.60"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("E26").Value = "  +7.17%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.66"
$ws.Range("E27").Value = "  +4.64%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").Value = "'9.66"
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").Value = "'740.28"
$ws.Range("E30").Value = "  +4.39%  "
$ws.Range("D31").Value = "'13.65"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").Value = "'0.131"
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("D34").Value = "'43.72"
$ws.Range("E34").Value = "  +8.87%  "
$ws.Range("E35").Value = "  +4.61%  "
$ws.Range("D36").Value = "'57.97"
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "'0.0₃0782"
$ws.Range("E38").Value = "  +14.21%  "
$ws.Range("D39").Value = "'5.35"
$ws.Range("E39").Value = "  -7.71%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "'3.25"
$ws.Range("E40").Value = "  +12.87%  "
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.336"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("E45").Value = "  +6.20%  "
$ws.Range("E46").Value = "  +4.99%  "
$ws.Range("D47").Value = "'3.43"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  -5.15%  "
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("D51").Value = "'143.57"
$ws.Range("E51").Value = "  +0.44%  "
